$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-08-22 Friday"; new = "2025-08-23 Saturday"},
    @{old = "407×5=2035"; new = "570×4=2280"},
    @{old = "811×3=2433"; new = "966×7=6762"},
    @{old = "550×8=4400"; new = "450×7=3150"},
    @{old = "213×6=1278"; new = "435×4=1740"},
    @{old = "918×9=8262"; new = "269×8=2152"},
    @{old = "274×4=1096"; new = "286×4=1144"},
    @{old = "673×2=1346"; new = "387×7=2709"},
    @{old = "708×8=5664"; new = "803×9=7227"},
    @{old = "364×8=2912"; new = "982×9=8838"},
    @{old = "312×7=2184"; new = "908×2=1816"},
    @{old = "836×6=5016"; new = "452×4=1808"},
    @{old = "477×9=4293"; new = "366×5=1830"},
    @{old = "755×6=4530"; new = "909×8=7272"},
    @{old = "114×2=228"; new = "710×8=5680"},
    @{old = "840×7=5880"; new = "921×7=6447"},
    @{old = "506×3=1518"; new = "120×9=1080"},
    @{old = "369×4=1476"; new = "640×8=5120"},
    @{old = "501×7=3507"; new = "647×4=2588"},
    @{old = "138×5=690"; new = "829×5=4145"},
    @{old = "319×4=1276"; new = "265×5=1325"},
    @{old = "760×9=6840"; new = "700×7=4900"},
    @{old = "926×3=2778"; new = "245×2=490"},
    @{old = "224×5=1120"; new = "281×9=2529"},
    @{old = "759×7=5313"; new = "418×7=2926"},
    @{old = "534×5=2670"; new = "350×2=700"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
